$d = $word.ActiveDocument

# --- Create the three new character styles ---
$ganStyle = $d.Styles.Add("GaNStyle", 2)
$ganStyle.Font.NameAscii = "Calibri"
$ganStyle.Font.Name = "Calibri"
$ganStyle.Font.Size = 14

$ganParagraph = $d.Styles.Add("GaNParagraph", 2)
$ganParagraph.Font.NameAscii = "Calibri"
$ganParagraph.Font.Name = "Calibri"
$ganParagraph.Font.Size = 10

$ganLinks = $d.Styles.Add("GaNLinks", 2)
$ganLinks.Font.NameAscii = "Calibri"
$ganLinks.Font.Name = "Calibri"
$ganLinks.Font.Bold = $true
$ganLinks.Font.Color = 8388608
$ganLinks.Font.Size = 9.5
$ganLinks.Font.Underline = 1

# --- Apply GaNStyle to every "2022 Campaign Dates..." run (4 occurrences) ---
$rng = $d.Content
$count = 0
while ($rng.Find.Execute(" 2022 Campaign Dates that use Leo constellation: April 14-23, May 14-23", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
    $count = $count + 1
    if ($count -gt 20) { break }
}

# --- Apply GaNParagraph to the "You are participating..." run ---
$rng2 = $d.Content
if ($rng2.Find.Execute("You are participating in a global campaign to observe and record the faintest stars visible as a means of measuring light pollution in a given location. By locating and observing the constellation Leo constellation in the night sky and comparing it to stellar charts, people from around the world will learn how the lights in their community contribute to light pollution. Your contributions to the online database will document the visible nighttime sky.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the link run ---
$rng3 = $d.Content
if ($rng3.Find.Execute("(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng3.Style = "GaNLinks"
}

Write-Output "Campaign date runs styled: $count"
